# This workbook is a pandas "date_proximity" merge-result fixture.
# The underlying refactor renamed the auxiliary merge-result columns
# (emitted by macpie's new core Dataset object) from:
#   _merge          -> _mp_merge
#   _diff_days      -> _mp_diff_days
#   _abs_diff_days  -> _mp_abs_diff_days
# and dropped the now-unused _duplicates boolean column entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing "_duplicates" column (M) - no longer produced.
$ws.Columns.Item(13).Delete()

# Rename the remaining merge-result headers to their "_mp_" prefixed names.
$ws.Range("J1").Value = "_mp_merge"
$ws.Range("K1").Value = "_mp_diff_days"
$ws.Range("L1").Value = "_mp_abs_diff_days"

# Let the now wider header text re-flow the column widths.
$ws.Columns.Item(10).ColumnWidth = 9.42
$ws.Columns.Item(11).ColumnWidth = 11.09
$ws.Columns.Item(12).ColumnWidth = 14.42
